$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteAll = -4104

# --- Column A labels shift up by one row (13..22 <- 14..23) ---
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A13").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A15").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A17").Copy() | Out-Null
$ws.Range("A16").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A17").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A20").Copy() | Out-Null
$ws.Range("A19").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A20").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A21").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A22").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("A23").ClearContents()

# --- B/C content that duplicates another existing cell, copied to preserve shared-string + style ---
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C14").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B10").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C10").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C18").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B16").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C17").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B23").PasteSpecial($xlPasteAll) | Out-Null
$ws.Range("C24").Copy() | Out-Null
$ws.Range("C23").PasteSpecial($xlPasteAll) | Out-Null

# --- Brand-new literal text values ---
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Clear cells that no longer have content ---
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()

# --- Row heights ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 15
$ws.Rows.Item(23).RowHeight = 30

# --- Remove the now-empty trailing row ---
$ws.Range("A24:C24").EntireRow.Delete()
